{"js": "// This script rewrites the Java stack-trace text in the second paragraph of the\n// document (produced by M2Doc's \"conditionRuntimeException\" test fixture).\n// The whole stack trace lives in a single run/<w:t> (bold, red). We locate that\n// block using two unique anchor strings (its first and last line) and replace\n// the entire block in one shot with the updated stack trace text that matches\n// the POI 3.16 -> 3.17 trace (updated line numbers, rebuilt JUnit/reflection\n// frames, and an extra nested Suite/RunAfters sequence), as described in the\n// commit (\"Fixed #253 Moving from POI 3.16 to 3.17.\").\n\nconst body = context.document.body;\n\n// Unique anchors marking the very start and very end of the stack-trace block.\nconst startResults = body.search(\"divOp(java.lang.Integer,java.lang.Integer) with arguments\", { matchCase: true });\nconst endResults = body.search(\"RemoteTestRunner.main(RemoteTestRunner.java:192)\", { matchCase: true });\nawait context.sync();\n\nif (startResults.items.length !== 1 || endResults.items.length !== 1) {\n  throw new Error(\"Expected to find exactly one start/end anchor for the stack trace block.\");\n}\n\nconst startRange = startResults.items[0];\nconst endRange = endResults.items[0];\nconst fullRange = startRange.expandTo(endRange);\n\nconst newStackTrace = \"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\\n\\t/ by zero\\njava.lang.ArithmeticException: / by zero\\n\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\\n\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\\n\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\\n\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\\n\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\\n\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:183)\\n\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\\n\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1267)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1)\\n\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\\n\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\\n\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\\n\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\\n\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\\n\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\\n\\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\\n\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\\n\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\\n\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\\n\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\\n\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\\n\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\\n\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\\n\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\\n\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\\n\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\\n\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\\n\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\\n\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\\n\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\\n\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\";\n\nfullRange.insertText(newStackTrace, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# This script rewrites the Java stack-trace text in the second paragraph of the\n# document (produced by M2Doc's \"conditionRuntimeException\" test fixture).\n# The whole stack trace lives in a single run (bold, red). We locate that block\n# using two unique anchor strings (its first and last line) and replace the\n# entire block in one shot with the updated stack trace text that matches the\n# POI 3.16 -> 3.17 trace (updated line numbers, rebuilt JUnit/reflection\n# frames, and an extra nested Suite/RunAfters sequence), as described in the\n# commit (\"Fixed #253 Moving from POI 3.16 to 3.17.\").\n\n$d = $word.ActiveDocument\n\n# Locate the start of the stack trace block.\n$startFind = $d.Range(0, $d.Content.End)\n$startFind.Find.Execute(\"divOp(java.lang.Integer,java.lang.Integer) with arguments\") | Out-Null\n$startPos = $startFind.Start\n\n# Locate the end of the stack trace block.\n$endFind = $d.Range(0, $d.Content.End)\n$endFind.Find.Execute(\"RemoteTestRunner.main(RemoteTestRunner.java:192)\") | Out-Null\n$endPos = $endFind.End\n\n# Range spanning the whole stack trace (keeps the run/paragraph formatting,\n# the trailing <w:br/> and the following paragraphs untouched).\n$rng = $d.Range($startPos, $endPos)\n\n$newStackTrace = @'\ndivOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\n\t/ by zero\njava.lang.ArithmeticException: / by zero\n\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\n\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:183)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\n\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1267)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\n\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\n\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\n\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\n\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\n\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\n\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\n'@\n\n$rng.Text = $newStackTrace\n"}
